$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column K (11), shifting K:P to L:Q
$ws.Columns.Item(11).Insert()

# Give the new column K the same width as column J, so Excel merges them
# into a single <col min="10" max="11" .../> element, matching a native
# "insert column" width-inheritance behavior.
$ws.Columns.Item(11).ColumnWidth = $ws.Columns.Item(10).ColumnWidth

# Apply the same style as column J (s="5") to the new K column cells
$ws.Range("J1:J13").Copy()
$ws.Range("K1:K13").PasteSpecial(-4122) # xlPasteFormats

# Set values in the exact order they first appear so the shared-strings
# table is built in the same sequence as the source edit.
$ws.Cells.Item(1, 11).Value = "azdhsId"
$ws.Cells.Item(11, 11).Value = "vacc35_1"
$ws.Cells.Item(3, 11).Value = "vacc34_1"
$ws.Cells.Item(4, 11).Value = "vacc34_1"
$ws.Cells.Item(9, 11).Value = "vacc26_1"
$ws.Cells.Item(10, 11).Value = "vacc22_1"
$ws.Cells.Item(2, 11).Value = "vacc18_1"
$ws.Cells.Item(7, 11).Value = "vacc15_1"
$ws.Cells.Item(8, 11).Value = "vacc15_1"
$ws.Cells.Item(6, 11).Value = "vacc10_1"
$ws.Cells.Item(5, 11).Value = "vacc8_1"
$ws.Cells.Item(12, 11).Value = "vacc38_1"
$ws.Cells.Item(13, 11).Value = "vacc39_1"

# Update selection to match target
$ws.Range("F21").Select()
